# edit.ps1 - applies the "crispian" diff via Word COM-interop
$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Paragraph 1: append two trailing spaces to the existing text,
#    then add three red-colored runs forming "(This is a change - Version for main branch)"
# ---------------------------------------------------------------
$para1 = $d.Paragraphs(1).Range
$para1Xml = '<w:p w14:paraId="5ADF5830" w14:textId="42E3A3E7" w:rsidR="00384372" w:rsidRDefault="00094D0B"><w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>(This is a change – Ve</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>rsion for main branch</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>)</w:t></w:r></w:p>'
$para1.InsertXML($para1Xml)

# ---------------------------------------------------------------
# 2) "Crispian's Day speech..." paragraph: restructure runs/proofErr
#    and merge the trailing " Henry V ... [Source - Wikipedia]" runs
#    into a single run.
# ---------------------------------------------------------------
$para4 = $d.Paragraphs(4).Range
$para4Xml = '<w:p w14:paraId="21745324" w14:textId="3B839C2B" w:rsidR="00347660" w:rsidRDefault="00AD0DA9" w:rsidP="00347660"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Crispian’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Day speech from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Shakespear’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Henry V [Source – Wikipedia]</w:t></w:r></w:p>'
$para4.InsertXML($para4Xml)

# ---------------------------------------------------------------
# 3) Add a new, empty paragraph styled "larger" right after the
#    final "...Saint Crispin's day." paragraph.
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$insertPoint = $d.Range($lastPara.End, $lastPara.End)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="larger"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0"/></w:pPr></w:p>'
$insertPoint.InsertXML($newParaXml)

# ---------------------------------------------------------------
# 4) styles.xml: mark "Normal (Web)" as semi-hidden, and remove the
#    now-unused "apple-converted-space" and "Hyperlink" character
#    styles. Delete by descending index so earlier deletions don't
#    shift the index of styles not yet removed.
# ---------------------------------------------------------------
$normalWeb = $d.Styles("Normal (Web)")
$normalWeb.Visibility = $false

$hyperlinkIdx = 0
$appleIdx = 0
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $styleName = $d.Styles($i).NameLocal
    if ($styleName -eq "Hyperlink") { $hyperlinkIdx = $i }
    if ($styleName -eq "apple-converted-space") { $appleIdx = $i }
}
if ($hyperlinkIdx -gt $appleIdx) {
    $d.Styles($hyperlinkIdx).Delete()
    $d.Styles($appleIdx).Delete()
} else {
    $d.Styles($appleIdx).Delete()
    $d.Styles($hyperlinkIdx).Delete()
}

Write-Host "edit applied"
